$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5954846143722534
$ws.Range("B1").Value = 0.9313179850578308
$ws.Range("C1").Value = 2.416607856750488
$ws.Range("D1").Value = 6.362005710601807
$ws.Range("E1").Value = 2.140755176544189
